$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.819.01"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "2.659.33"
$ws.Range("E3").Value = "  +2.57%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'600.98"
$ws.Range("E5").Value = "  +2.06%  "

$ws.Range("D6").Value = "'155.53"
$ws.Range("E6").Value = "  +4.20%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D9").Value = "2.659.32"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").Value = "'0.139"
$ws.Range("E10").Value = "  +12.67%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "'5.24"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("D14").Value = "'28.08"
$ws.Range("E14").Value = "  +3.47%  "

$ws.Range("E15").Value = "  +6.32%  "

$ws.Range("D16").Value = "3.142.59"
$ws.Range("E16").Value = "  +2.62%  "

$ws.Range("D17").Value = "68.697.35"
$ws.Range("E17").Value = "  +2.63%  "

$ws.Range("D18").Value = "2.658.01"
$ws.Range("E18").Value = "  +2.50%  "

$ws.Range("E19").Value = "  +4.46%  "

$ws.Range("D20").Value = "'367.62"
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("D21").Value = "'7.46"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E23").Value = "  +1.38%  "

$ws.Range("E24").Value = "  +5.36%  "

$ws.Range("D25").Value = "'72.72"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  +1.33%  "

$ws.Range("E28").Value = "  +8.82%  "

$ws.Range("D29").Value = "2.787.97"

$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").Value = "'577.58"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "  +5.00%  "

$ws.Range("D33").Value = "'7.99"
$ws.Range("E33").Value = "  +5.30%  "

$ws.Range("E34").Value = "  +3.62%  "

$ws.Range("E35").Value = "  +5.37%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "'1.55"
$ws.Range("E37").Value = "  +4.22%  "

$ws.Range("D38").Value = "'158.74"
$ws.Range("E38").Value = "  +1.48%  "

$ws.Range("E39").Value = "  +5.62%  "

$ws.Range("D40").Value = "'19.30"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("E41").Value = "  +4.68%  "

$ws.Range("D42").Value = "'0.369"
$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("D43").Value = "'2.67"
$ws.Range("E43").Value = "  +7.93%  "

$ws.Range("D44").Value = "'17.76"
$ws.Range("E44").Value = "  +5.71%  "

$ws.Range("D45").Value = "0.0₆0321"
$ws.Range("E45").Value = "  +13.66%  "

$ws.Range("D46").Value = "'40.69"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("D48").Value = "'157.05"
$ws.Range("E48").Value = "  +3.53%  "

$ws.Range("D49").Value = "'3.75"
$ws.Range("E49").Value = "  +1.13%  "

$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("D51").Value = "'22.05"
$ws.Range("E51").Value = "  +4.01%  "
